# Update the latest referee / linesman KHL stats snapshot
# (refreshed Games/PIM figures + as_of_utc timestamps) on the
# "Главные" (referees) and "Линейные" (linesmen) sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: Главные ---
$ws = $wb.Worksheets.Item("Главные")
# Row 2: Akuzovskiy Nikolay
$ws.Range("AA2").Value = "2025-12-02 03:04:09"
# Row 3: Belov Aleksey
$ws.Range("C3").Value = 30
$ws.Range("D3").Value = 521
$ws.Range("E3").Value = 235
$ws.Range("G3").Value = 17.37
$ws.Range("H3").Value = 7.83
$ws.Range("I3").Value = 9.529999999999999
$ws.Range("J3").Value = 115
$ws.Range("AA3").Value = "2025-12-02 03:04:09"
# Row 4: Belyaev Dmitriy
$ws.Range("AA4").Value = "2025-12-02 03:04:09"
# Row 5: Belyaev Sergey
$ws.Range("C5").Value = 30
$ws.Range("D5").Value = 483
$ws.Range("E5").Value = 244
$ws.Range("F5").Value = 239
$ws.Range("G5").Value = 16.1
$ws.Range("H5").Value = 8.130000000000001
$ws.Range("J5").Value = 117
$ws.Range("K5").Value = 112
$ws.Range("AA5").Value = "2025-12-02 03:04:09"
# Row 6: Birin Viktor
$ws.Range("AA6").Value = "2025-12-02 03:04:09"
# Row 7: Vasilev Aleksey
$ws.Range("AA7").Value = "2025-12-02 03:04:09"
# Row 8: Gamaley Evgeniy
$ws.Range("AA8").Value = "2025-12-02 03:04:09"
# Row 9: Gashilov Viktor
$ws.Range("AA9").Value = "2025-12-02 03:04:09"
# Row 10: Gofman Anton
$ws.Range("C10").Value = 20
$ws.Range("D10").Value = 338
$ws.Range("E10").Value = 172
$ws.Range("F10").Value = 166
$ws.Range("G10").Value = 16.9
$ws.Range("H10").Value = 8.6
$ws.Range("I10").Value = 8.300000000000001
$ws.Range("J10").Value = 86
$ws.Range("K10").Value = 73
$ws.Range("AA10").Value = "2025-12-02 03:04:09"
# Row 11: Dudarov Aleksandr
$ws.Range("AA11").Value = "2025-12-02 03:04:09"
# Row 12: Kochetov Evgeniy
$ws.Range("C12").Value = 20
$ws.Range("D12").Value = 331
$ws.Range("E12").Value = 143
$ws.Range("G12").Value = 16.55
$ws.Range("H12").Value = 7.15
$ws.Range("I12").Value = 9.4
$ws.Range("J12").Value = 59
$ws.Range("AA12").Value = "2025-12-02 03:04:09"
# Row 13: Krasotin Nikolay
$ws.Range("C13").Value = 11
$ws.Range("D13").Value = 167
$ws.Range("E13").Value = 94
$ws.Range("F13").Value = 73
$ws.Range("G13").Value = 15.18
$ws.Range("H13").Value = 8.550000000000001
$ws.Range("I13").Value = 6.64
$ws.Range("J13").Value = 47
$ws.Range("K13").Value = 34
$ws.Range("AA13").Value = "2025-12-02 03:04:09"
# Row 14: Lavrentev Anton
$ws.Range("AA14").Value = "2025-12-02 03:04:09"
# Row 15: Lazarev Gleb
$ws.Range("AA15").Value = "2025-12-02 03:04:09"
# Row 16: Morozov Sergey
$ws.Range("AA16").Value = "2025-12-02 03:04:09"
# Row 17: Mochalov Vladimir
$ws.Range("AA17").Value = "2025-12-02 03:04:09"
# Row 18: Naumov Denis
$ws.Range("AA18").Value = "2025-12-02 03:04:09"
# Row 19: Ovchinnikov Pavel
$ws.Range("AA19").Value = "2025-12-02 03:04:09"
# Row 20: Oskirko Yuriy
$ws.Range("C20").Value = 30
$ws.Range("D20").Value = 517
$ws.Range("E20").Value = 226
$ws.Range("F20").Value = 291
$ws.Range("G20").Value = 17.23
$ws.Range("H20").Value = 7.53
$ws.Range("I20").Value = 9.699999999999999
$ws.Range("J20").Value = 103
$ws.Range("K20").Value = 108
$ws.Range("AA20").Value = "2025-12-02 03:04:09"
# Row 21: Romasko Evgeniy
$ws.Range("AA21").Value = "2025-12-02 03:04:09"
# Row 22: Svetilov Aleksey
$ws.Range("AA22").Value = "2025-12-02 03:04:09"
# Row 23: Sergeev Aleksandr V.
$ws.Range("C23").Value = 19
$ws.Range("D23").Value = 240
$ws.Range("E23").Value = 95
$ws.Range("F23").Value = 145
$ws.Range("G23").Value = 12.63
$ws.Range("H23").Value = 5
$ws.Range("I23").Value = 7.63
$ws.Range("J23").Value = 45
$ws.Range("K23").Value = 60
$ws.Range("AA23").Value = "2025-12-02 03:04:09"
# Row 24: Sidorenko Maksim
$ws.Range("AA24").Value = "2025-12-02 03:04:09"
# Row 25: Soin Aleksandr
$ws.Range("C25").Value = 30
$ws.Range("D25").Value = 475
$ws.Range("E25").Value = 230
$ws.Range("F25").Value = 245
$ws.Range("G25").Value = 15.83
$ws.Range("H25").Value = 7.67
$ws.Range("J25").Value = 110
$ws.Range("K25").Value = 115
$ws.Range("AA25").Value = "2025-12-02 03:04:09"
# Row 26: Spirin Viktor
$ws.Range("AA26").Value = "2025-12-02 03:04:09"

# --- Sheet: Линейные ---
$ws = $wb.Worksheets.Item("Линейные")
# Row 2: Baranov Nikita
$ws.Range("AA2").Value = "2025-12-02 03:04:09"
# Row 3: Bersenyov Maksim
$ws.Range("AA3").Value = "2025-12-02 03:04:09"
# Row 4: Boldyrev Sergey
$ws.Range("C4").Value = 14
$ws.Range("D4").Value = 208
$ws.Range("E4").Value = 94
$ws.Range("F4").Value = 114
$ws.Range("G4").Value = 14.86
$ws.Range("H4").Value = 6.71
$ws.Range("I4").Value = 8.140000000000001
$ws.Range("J4").Value = 47
$ws.Range("K4").Value = 47
$ws.Range("AA4").Value = "2025-12-02 03:04:09"
# Row 5: Bulychev Egor
$ws.Range("C5").Value = 16
$ws.Range("D5").Value = 226
$ws.Range("E5").Value = 118
$ws.Range("F5").Value = 108
$ws.Range("G5").Value = 14.13
$ws.Range("H5").Value = 7.38
$ws.Range("I5").Value = 6.75
$ws.Range("J5").Value = 59
$ws.Range("K5").Value = 54
$ws.Range("AA5").Value = "2025-12-02 03:04:09"
# Row 6: Buturlin Vladimir
$ws.Range("C6").Value = 18
$ws.Range("D6").Value = 309
$ws.Range("E6").Value = 145
$ws.Range("F6").Value = 164
$ws.Range("G6").Value = 17.17
$ws.Range("H6").Value = 8.06
$ws.Range("I6").Value = 9.109999999999999
$ws.Range("J6").Value = 65
$ws.Range("K6").Value = 77
$ws.Range("AA6").Value = "2025-12-02 03:04:09"
# Row 7: Bukharov Nikita
$ws.Range("AA7").Value = "2025-12-02 03:04:09"
# Row 8: Vilyugin Nikita
$ws.Range("AA8").Value = "2025-12-02 03:04:09"
# Row 9: Golovlyov Dmitriy
$ws.Range("AA9").Value = "2025-12-02 03:04:09"
# Row 10: Gribovskiy Nikita
$ws.Range("AA10").Value = "2025-12-02 03:04:09"
# Row 11: Egorov Sergey
$ws.Range("AA11").Value = "2025-12-02 03:04:09"
# Row 12: Zaytsev Valentin
$ws.Range("AA12").Value = "2025-12-02 03:04:09"
# Row 13: Ivanichkin Ivan
$ws.Range("AA13").Value = "2025-12-02 03:04:09"
# Row 14: Ivanov Yuriy
$ws.Range("C14").Value = 28
$ws.Range("D14").Value = 448
$ws.Range("E14").Value = 227
$ws.Range("G14").Value = 16
$ws.Range("H14").Value = 8.109999999999999
$ws.Range("I14").Value = 7.89
$ws.Range("J14").Value = 111
$ws.Range("AA14").Value = "2025-12-02 03:04:09"
# Row 15: Kolyasnikov Ilya
$ws.Range("AA15").Value = "2025-12-02 03:04:09"
# Row 16: Kupriyanov Maksim
$ws.Range("C16").Value = 29
$ws.Range("D16").Value = 478
$ws.Range("E16").Value = 220
$ws.Range("G16").Value = 16.48
$ws.Range("H16").Value = 7.59
$ws.Range("I16").Value = 8.9
$ws.Range("J16").Value = 100
$ws.Range("AA16").Value = "2025-12-02 03:04:09"
# Row 17: Litvinov Aleksandr
$ws.Range("AA17").Value = "2025-12-02 03:04:09"
# Row 18: Novikov Nikita
$ws.Range("C18").Value = 32
$ws.Range("D18").Value = 529
$ws.Range("E18").Value = 247
$ws.Range("F18").Value = 282
$ws.Range("G18").Value = 16.53
$ws.Range("H18").Value = 7.72
$ws.Range("I18").Value = 8.81
$ws.Range("J18").Value = 116
$ws.Range("K18").Value = 121
$ws.Range("AA18").Value = "2025-12-02 03:04:09"
# Row 19: Parikov Yaroslav
$ws.Range("C19").Value = 28
$ws.Range("D19").Value = 459
$ws.Range("E19").Value = 220
$ws.Range("F19").Value = 239
$ws.Range("G19").Value = 16.39
$ws.Range("H19").Value = 7.86
$ws.Range("I19").Value = 8.539999999999999
$ws.Range("J19").Value = 105
$ws.Range("K19").Value = 107
$ws.Range("AA19").Value = "2025-12-02 03:04:09"
# Row 20: Polyakov Nikita A.
$ws.Range("AA20").Value = "2025-12-02 03:04:09"
# Row 21: Sadovnikov Aleksandr
$ws.Range("AA21").Value = "2025-12-02 03:04:09"
# Row 22: Sedov Egor
$ws.Range("AA22").Value = "2025-12-02 03:04:09"
# Row 23: Serdyuk Nikita
$ws.Range("AA23").Value = "2025-12-02 03:04:09"
# Row 24: Sivov Dmitriy
$ws.Range("AA24").Value = "2025-12-02 03:04:09"
# Row 25: Skugarev Miroslav
$ws.Range("AA25").Value = "2025-12-02 03:04:09"
# Row 26: Slavikovskiy Roman
$ws.Range("AA26").Value = "2025-12-02 03:04:09"

